$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------
# 1) "Leather " + "Spaulder" (two runs around a spell-check proofErr
#    pair) -> single clean run "Leather Spaulder" with no proofing
#    markup left behind.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Leather Spaulder") {
        $r = $p.Range
        $r.Delete()
        $xml = "<w:p xmlns:w='$wNs'><w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr><w:r><w:t>Leather Spaulder</w:t></w:r></w:p>"
        $r.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------
# 2) Append a new "Skills:" section (Heading2) with descriptions for
#    each skill after the "+2 Political" paragraph at the end of the
#    document.
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "+2 Political") {
        $target = $p
    }
}

$insertPos = $target.Range.End
$r = $d.Range($insertPos, $insertPos)

$xml = @"
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='Heading2'/></w:pPr>
  <w:r><w:t>Skills</w:t></w:r>
  <w:r><w:t>:</w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
  <w:r><w:t>Athletics: Used for vaulting over rough terrain.</w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
  <w:r><w:t>Melee: Used to increase your chance to hit enemies with melee attacks.</w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
  <w:r><w:t>Ranged: Used to increase your accuracy with ranged weapons (not yet in game)</w:t></w:r>
  <w:bookmarkStart w:id='0' w:name='_GoBack'/>
  <w:bookmarkEnd w:id='0'/>
  <w:r><w:t>.</w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
  <w:r><w:t>Stealth: Used to hide from enemies better when using stealth.</w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
  <w:r><w:t>Mechanical: Used to build traps and turrets.</w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
  <w:r><w:t>Medicinal: Used to improve efficiency of med kits (not yet in game).</w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
  <w:r><w:t>Historical: Used for more knowledge for missions and during dialogue (not yet in game).</w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
</w:p>
<w:p xmlns:w='$wNs'>
  <w:pPr><w:pStyle w:val='NoSpacing'/></w:pPr>
  <w:r><w:t>Political: Used for persuasive skills such as Orator features and dialogue responses.</w:t></w:r>
</w:p>
"@

$r.InsertXML($xml)
